$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.089.32"
$ws.Range("E2").Value = "  -0.76%  "
$ws.Range("D3").Value = "1.651.76"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  -0.42%  "
$ws.Range("D5").Value = "'218.12"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").Value = "'0.5212"
$ws.Range("E6").Value = "  -2.17%  "
$ws.Range("E7").Value = "  -0.40%  "
$ws.Range("D8").Value = "'0.2622"
$ws.Range("E8").Value = "  -0.99%  "
$ws.Range("D9").Value = "'0.06297"
$ws.Range("E9").Value = "  -1.28%  "
$ws.Range("D10").Value = "'20.52"
$ws.Range("E10").Value = "  -0.16%  "
$ws.Range("D11").Value = "'0.07803"
$ws.Range("D12").Value = "'4.483"
$ws.Range("E12").Value = "  -1.83%  "
$ws.Range("D13").Value = "1.650.13"
$ws.Range("E13").Value = "  -0.89%  "
$ws.Range("D14").Value = "1.877.92"
$ws.Range("E14").Value = "  -0.76%  "
$ws.Range("D15").Value = "'0.5555"
$ws.Range("E15").Value = "  +0.55%  "
$ws.Range("D16").Value = "0.0₅8010"
$ws.Range("E16").Value = "  -2.38%  "
$ws.Range("D17").Value = "'64.90"
$ws.Range("E17").Value = "  -1.11%  "
$ws.Range("D18").Value = "26.077.19"
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("D20").Value = "'4.643"
$ws.Range("E20").Value = "  -0.98%  "
$ws.Range("D21").Value = "'194.78"
$ws.Range("E21").Value = "  +0.72%  "
$ws.Range("D22").Value = "'10.09"
$ws.Range("E22").Value = "  -1.14%  "
$ws.Range("D23").Value = "'5.952"
$ws.Range("E23").Value = "  -1.34%  "
$ws.Range("D24").Value = "'1.006"
$ws.Range("E24").Value = "  -0.42%  "
$ws.Range("D25").Value = "'146.71"
$ws.Range("D26").Value = "'0.1206"
$ws.Range("E26").Value = "  -2.04%  "
$ws.Range("D27").Value = "'7.185"
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("D28").Value = "'15.91"
$ws.Range("E28").Value = "  -1.30%  "
$ws.Range("D29").Value = "'1.470"
$ws.Range("E29").Value = "  -0.68%  "
$ws.Range("D30").Value = "'0.05709"
$ws.Range("E30").Value = "  -2.52%  "
$ws.Range("D31").Value = "'1.266"
$ws.Range("E31").Value = "  -1.18%  "
$ws.Range("D32").Value = "'3.487"
$ws.Range("E32").Value = "  -3.34%  "
$ws.Range("D33").Value = "'3.365"
$ws.Range("E33").Value = "  +2.63%  "
$ws.Range("D34").Value = "'1.592"
$ws.Range("E34").Value = "  -1.02%  "
$ws.Range("D35").Value = "'2.799"
$ws.Range("E35").Value = "  -0.88%  "
$ws.Range("D36").Value = "'0.9502"
$ws.Range("E36").Value = "  -1.22%  "
$ws.Range("E37").Value = "  -0.32%  "
$ws.Range("D38").Value = "'0.5661"
$ws.Range("E38").Value = "  -2.40%  "
$ws.Range("D39").Value = "'0.01588"
$ws.Range("E39").Value = "  -1.22%  "
$ws.Range("D40").Value = "'5.971"
$ws.Range("E40").Value = "  +2.55%  "
$ws.Range("D41").Value = "1.057.47"
$ws.Range("E41").Value = "  +0.71%  "
$ws.Range("D42").Value = "'1.005"
$ws.Range("E42").Value = "  -0.42%  "
$ws.Range("D43").Value = "'0.8422"
$ws.Range("E43").Value = "  -2.67%  "
$ws.Range("D44").Value = "'103.65"
$ws.Range("E44").Value = "  -0.85%  "
$ws.Range("D45").Value = "1.789.73"
$ws.Range("E45").Value = "  -0.73%  "
$ws.Range("D46").Value = "'57.39"
$ws.Range("E46").Value = "  -0.53%  "
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").Value = "'1.005"
$ws.Range("E47").Value = "  -0.35%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.05331"
$ws.Range("E48").Value = "  +3.22%  "
$ws.Range("D49").Value = "'0.4354"
$ws.Range("E49").Value = "  -0.64%  "
$ws.Range("E50").Value = "  -3.09%  "
$ws.Range("D51").Value = "'7.950"
$ws.Range("E51").Value = "  -0.98%  "
